$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (Q0)
$ws.Range("B2").Value = 0.06284567650902455
$ws.Range("C2").Value = 0.6441637244774302
$ws.Range("D2").Value = 0.7722581226213636
$ws.Range("E2").Value = 0.8787821815565924
$ws.Range("F2").Value = 0.8852540394919396
$ws.Range("G2").Value = 51

# Row 3 (Q1)
$ws.Range("B3").Value = -0.1148980440451865
$ws.Range("C3").Value = 0.812265897056845
$ws.Range("D3").Value = 1.412039409307321
$ws.Range("E3").Value = 1.188292644640756
$ws.Range("F3").Value = 1.194732421917788
$ws.Range("G3").Value = 50

# Row 4 (Q2)
$ws.Range("B4").Value = 0.06179048950849853
$ws.Range("C4").Value = 0.7574901860148398
$ws.Range("D4").Value = 0.970307845326831
$ws.Range("E4").Value = 0.9850420525677221
$ws.Range("F4").Value = 0.9932899800654412
$ws.Range("G4").Value = 49

# Row 5 (Q3)
$ws.Range("B5").Value = -0.1253060257331123
$ws.Range("C5").Value = 0.8314090100821168
$ws.Range("D5").Value = 1.348981947049162
$ws.Range("E5").Value = 1.161456821000747
$ws.Range("F5").Value = 1.166896745183917
$ws.Range("G5").Value = 48

# Row 6 (Q4)
$ws.Range("B6").Value = 0.02745591576579163
$ws.Range("C6").Value = 0.7011270665634526
$ws.Range("D6").Value = 0.9424265274435161
$ws.Range("E6").Value = 0.9707865509181285
$ws.Range("F6").Value = 0.9808893137291498
$ws.Range("G6").Value = 47

# Row 7 (Q5)
$ws.Range("B7").Value = -0.1781507454251396
$ws.Range("C7").Value = 0.7419848673622177
$ws.Range("D7").Value = 1.093258070565441
$ws.Range("E7").Value = 1.045589819463369
$ws.Range("F7").Value = 1.04168600081914
$ws.Range("G7").Value = 46

# Row 8 (Q6)
$ws.Range("B8").Value = -0.03744529641089314
$ws.Range("C8").Value = 0.7201853156618916
$ws.Range("D8").Value = 1.085640148578663
$ws.Range("E8").Value = 1.041940568640392
$ws.Range("F8").Value = 1.053033603948734
$ws.Range("G8").Value = 45

# Row 9 (Q7)
$ws.Range("B9").Value = -0.1326544656744229
$ws.Range("C9").Value = 0.8526928699447087
$ws.Range("D9").Value = 1.204565724579542
$ws.Range("E9").Value = 1.097527095145966
$ws.Range("F9").Value = 1.102076420364545
$ws.Range("G9").Value = 44

# Row 10 (Q8)
$ws.Range("B10").Value = -0.06436346714910372
$ws.Range("C10").Value = 0.8166505515847911
$ws.Range("D10").Value = 1.166152122974641
$ws.Range("E10").Value = 1.07988523602031
$ws.Range("F10").Value = 1.090722860833269
$ws.Range("G10").Value = 43

# Row 11 (new row, Q9) - copy formatting from row 10's label cell, then set values
$ws.Range("A10").Copy()
$ws.Range("A11").PasteSpecial(-4122, -4142, $false, $false)

$ws.Range("A11").Value = "Q9"
$ws.Range("B11").Value = -0.1139447951980307
$ws.Range("C11").Value = 0.7616368124151508
$ws.Range("D11").Value = 0.970461329011336
$ws.Range("E11").Value = 0.9851199566607794
$ws.Range("F11").Value = 0.9903691395028201
$ws.Range("G11").Value = 42

Write-Host "done"
